$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2: value 11 -> 10.5, and restyle to match C5 (style index 3)
$ws.Range("C2").Value = 10.5
$ws.Range("C2").Style = $ws.Range("C5").Style

# Update C3: value 10 -> 9, and restyle to match C5 (style index 3)
$ws.Range("C3").Value = 9
$ws.Range("C3").Style = $ws.Range("C5").Style

# Update the selected cell to F8
$ws.Range("F8").Select()
